# Append the 10/28/2025 profit record as a new row (72) at the bottom of
# the daily Date/Profit log on Sheet1, per "Update profit files after
# running on 2025-10-28".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (e.g. "10/27/2025" in A71), not a
# real Excel date serial. Force the new cell to Text format before writing
# the string so COM doesn't auto-convert "10/28/2025" into a date value,
# then drop the format stamp again so the cell ends up plain/unstyled like
# its neighbours.
$dateCell = $ws.Range("A72")
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/28/2025"
$dateCell.ClearFormats()

# Column B holds the plain numeric profit figure for that day.
$ws.Range("B72").Value = 11816.54
